$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "Arial 10pt black" direct formatting (currently applied to
# E13) by copying it onto the rows that will hold the new data (E2:E5) before the
# old rows are removed.
$ws.Range("E13").Copy()
$ws.Range("E2:E5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Drop the old data rows 6-13; rows 2-5 remain (with their original formatting)
# and get overwritten with the new values below.
$ws.Rows("6:13").Delete()

# New table body - write column C (datatype) before column B (attribute name) so
# that new shared-string entries are interned in the same order as the source
# workbook.
$ws.Cells.Item(2, 3).Value = "mdex:double"
$ws.Cells.Item(2, 2).Value = "TAX_AMT_PEMEX"
$ws.Cells.Item(2, 1).Value = 516
$ws.Cells.Item(2, 4).Value = 6
$ws.Cells.Item(2, 5).Value = "IVA en Pesos"

$ws.Cells.Item(3, 2).Value = "TAX_VALUE_PEMEX"
$ws.Cells.Item(3, 1).Value = 516
$ws.Cells.Item(3, 3).Value = "mdex:double"
$ws.Cells.Item(3, 4).Value = 6
$ws.Cells.Item(3, 5).Value = "IVA en DLS"

$ws.Cells.Item(4, 2).Value = "TAX_AMOUNT_PEMEX"
$ws.Cells.Item(4, 1).Value = 516
$ws.Cells.Item(4, 3).Value = "mdex:double"
$ws.Cells.Item(4, 4).Value = 6
$ws.Cells.Item(4, 5).Value = "IVA en DLS"

$ws.Cells.Item(5, 2).Value = "INVOICE_AMOUNT_WITH_TAX"
$ws.Cells.Item(5, 1).Value = 516
$ws.Cells.Item(5, 3).Value = "mdex:double"
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = "Monto Total en DLS"

# Column widths (A widens to fit the longer header text; B-E unchanged).
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(3).ColumnWidth = 15.5
$ws.Columns.Item(4).ColumnWidth = 9.833333333333334
$ws.Columns.Item(5).ColumnWidth = 32.666666666666664

# Match the final selection/active cell.
$ws.Range("E5").Select()
